# Sync file from Google Drive
# Updates ETA timestamps / minutes-to-arrival / a couple of flag & type
# fields on the three NextBus sheets, and appends a new "NextBus3" data
# row (row 15) to the NextBus2 and NextBus3 sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: write a value into a cell's column-O "MinutesToArrival" field
# while keeping it a genuine text value (the source data stores these
# as strings, e.g. "11", not numbers) by using Excel's leading-apostrophe
# text-entry convention - otherwise a numeric-looking string like "11"
# gets auto-coerced to a Number by the Value setter.
# ---------------------------------------------------------------------
function Set-TextValue($ws, [string]$addr, [string]$text) {
    $ws.Range($addr).Value = "'" + $text
}

# =======================================================================
# Sheet "NextBus1"
# =======================================================================
$ws1 = $wb.Worksheets.Item("NextBus1")

$ws1.Range("F2").Value = 45688.68462962963
Set-TextValue $ws1 "O2" "11"

$ws1.Range("F4").Value = 45688.68417824074
Set-TextValue $ws1 "O4" "11"

$ws1.Range("F5").Value = 45688.69063657407
Set-TextValue $ws1 "O5" "20"

$ws1.Range("L6").Value = "SD"

$ws1.Range("J7").Value = 1

$ws1.Range("F8").Value = 45688.68344907407
Set-TextValue $ws1 "O8" "10"

$ws1.Range("F10").Value = 45688.685
Set-TextValue $ws1 "O10" "12"

$ws1.Range("F11").Value = 45688.68221064815
Set-TextValue $ws1 "O11" "8"

$ws1.Range("F12").Value = 45688.68560185185
Set-TextValue $ws1 "O12" "13"

$ws1.Range("J15").Value = 1

# =======================================================================
# Sheet "NextBus2"
# =======================================================================
$ws2 = $wb.Worksheets.Item("NextBus2")

$ws2.Range("F2").Value = 45688.69228009259
Set-TextValue $ws2 "O2" "22"

$ws2.Range("F4").Value = 45688.68792824074
Set-TextValue $ws2 "O4" "16"

$ws2.Range("F5").Value = 45688.6971875
Set-TextValue $ws2 "O5" "29"

$ws2.Range("L6").Value = "DD"

$ws2.Range("F8").Value = 45688.69637731482
Set-TextValue $ws2 "O8" "28"

$ws2.Range("F10").Value = 45688.69920138889
Set-TextValue $ws2 "O10" "32"

$ws2.Range("F11").Value = 45688.68884259259
$ws2.Range("J11").Value = 1
Set-TextValue $ws2 "O11" "17"

$ws2.Range("F12").Value = 45688.69226851852
$ws2.Range("J12").Value = 1
Set-TextValue $ws2 "O12" "22"

# New row 15 - another NextBus3 departure appended to the feed
$ws2.Range("A15").Value = "NextBus3"
$ws2.Range("B15").Value = 75
$ws2.Range("C15").Value = 10009
$ws2.Range("D15").Value = "Bt Merah Int"
$ws2.Range("E15").Value = "SMRT"
$ws2.Range("G15").Value = 10009
$ws2.Range("H15").Value = "WAB"
$ws2.Range("I15").Value = "SEA"
$ws2.Range("J15").Value = 0
$ws2.Range("K15").Value = 44989
$ws2.Range("L15").Value = "SD"
$ws2.Range("M15").Value = 12109
$ws2.Range("N15").Value = "Opp Ngee Ann Poly"

# =======================================================================
# Sheet "NextBus3"
# =======================================================================
$ws3 = $wb.Worksheets.Item("NextBus3")

$ws3.Range("F2").Value = 45688.69228009259
Set-TextValue $ws3 "O2" "22"

$ws3.Range("F4").Value = 45688.68792824074
Set-TextValue $ws3 "O4" "16"

$ws3.Range("F5").Value = 45688.6971875
Set-TextValue $ws3 "O5" "29"

$ws3.Range("L6").Value = "DD"

$ws3.Range("F8").Value = 45688.69637731482
Set-TextValue $ws3 "O8" "28"

$ws3.Range("F10").Value = 45688.69920138889
Set-TextValue $ws3 "O10" "32"

$ws3.Range("F11").Value = 45688.68884259259
$ws3.Range("J11").Value = 1
Set-TextValue $ws3 "O11" "17"

$ws3.Range("F12").Value = 45688.69226851852
$ws3.Range("J12").Value = 1
Set-TextValue $ws3 "O12" "22"

# New row 15 - another NextBus3 departure appended to the feed
$ws3.Range("A15").Value = "NextBus3"
$ws3.Range("B15").Value = 75
$ws3.Range("C15").Value = 10009
$ws3.Range("D15").Value = "Bt Merah Int"
$ws3.Range("E15").Value = "SMRT"
$ws3.Range("G15").Value = 10009
$ws3.Range("H15").Value = "WAB"
$ws3.Range("I15").Value = "SEA"
$ws3.Range("J15").Value = 0
$ws3.Range("K15").Value = 44989
$ws3.Range("L15").Value = "SD"
$ws3.Range("M15").Value = 12109
$ws3.Range("N15").Value = "Opp Ngee Ann Poly"
